$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: add "% White" column header in I5:I6 (merged, matching H5:H6 formatting) ---
$ws.Range("H5:H6").Copy()
$ws.Range("I5:I6").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("I5").Value2 = "% White"
$ws.Range("I5:I6").Merge()

# --- Data rows 7-98: add "% White" = B / H formula ---
$ws.Range("I7:I98").Formula = "=B7/H7"
$ws.Range("I7:I98").Style = "Percent"

# --- Row 99 (Total row): add B99 = SUM(B7:B98) and I99 = B99/H99 ---
$ws.Range("B99").Formula = "=SUM(B7:B98)"
$ws.Range("B99").Style = $ws.Range("H99").Style
$ws.Range("B99").NumberFormat = $ws.Range("H99").NumberFormat
$ws.Range("I99").Formula = "=B99/H99"
$ws.Range("I99").Style = "Percent"

# --- Row 100: empty I100 cell carrying the Percent style ---
$ws.Range("I100").Style = "Percent"

# --- Update selection to match the target view ---
$ws.Range("M83").Select()
